# Fix(all scripts): Updating routes for __tmp__ files
# Rebuild the Huff-model input table: new supermarket columns (Hipermercado
# Metro Independencia, Tottus Mega Plaza, Tottus Los Olivos, Makro Plaza Lima
# Norte, Makro Comas) and new store rows (Z651, Z423A/Z423B split, Z398/Z396
# kept) with binary (0/1) flags replacing the old 0.4/0.1 weights.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column headers (row 1) ----
$ws.Range("B1").Value = "Supermercado"
$ws.Range("C1").Value = "Hipermercado Metro Independencia"
$ws.Range("D1").Value = "Plaza Vea Izaguirre"
$ws.Range("E1").Value = "Tottus Mega Plaza"
$ws.Range("F1").Value = "Plaza Vea Los Olivos"
$ws.Range("G1").Value = "Tottus Los Olivos"
$ws.Range("H1").Value = "Makro Plaza Lima Norte"
$ws.Range("I1").Value = "Makro Comas"

# ---- Row labels (column A) and data rows (B:I) ----
# Row 2 - Z408
$ws.Range("A2").Value = "Z408"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 0

# Row 3 - Z651
$ws.Range("A3").Value = "Z651"
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0

# Row 4 - Z407
$ws.Range("A4").Value = "Z407"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 0

# Row 5 - Z414
$ws.Range("A5").Value = "Z414"
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 0

# Row 6 - Z409
$ws.Range("A6").Value = "Z409"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 0

# Row 7 - Z403
$ws.Range("A7").Value = "Z403"
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 0

# Row 8 - Z412
$ws.Range("A8").Value = "Z412"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 0

# Row 9 - Z405
$ws.Range("A9").Value = "Z405"
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 0

# Row 10 - Z417
$ws.Range("A10").Value = "Z417"
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 0

# Row 11 - Z423B
$ws.Range("A11").Value = "Z423B"
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0

# Row 12 - Z411
$ws.Range("A12").Value = "Z411"
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0

# Row 13 - Z399
$ws.Range("A13").Value = "Z399"
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 1
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 0

# Row 14 - Z423A
$ws.Range("A14").Value = "Z423A"
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0

# Row 15 - Z398
$ws.Range("A15").Value = "Z398"
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = 1

# Row 16 - Z396 (new row)
$ws.Range("A16").Value = "Z396"
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 1

# ---- Update selection to match new extent ----
$ws.Range("A1:I16").Select()
